# Populate the worksheet with the new cell contents.
# Write order chosen so the shared-string table comes out in the same
# index order as the target workbook (this=0, Hello=1, this1=2).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "this"
$ws.Range("A1").Value = "Hello"
$ws.Range("B2").Value = "this1"

# Add the "Button 7" form-control button (macro trigger) roughly where the
# original lived: anchored from column G (idx 6) row 4 (idx 3) to column J
# (idx 9) row 9 (idx 8), converted from EMU to points (1 pt = 12700 EMU).
$shp = $ws.Shapes.AddFormControl(0, 303.75, 59.25, 144, 66)
$shp.Name = "Button 7"
$shp.OnAction = "RunPythonScript"
$shp.TextFrame.Characters().Text = "Button 7"

# Restore the author's last selection.
$ws.Range("I9").Select() | Out-Null
